$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---

# Overview sheet: status columns for zh-cn (E2) and de-de (F2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column (C2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column (C2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Adjust column widths to reflect the narrower status text ---
# (target raw width 13.4101845877511; closest value this host's column-width
# grid can express is 13.3333.., reached via ColumnWidth = 12.5)

# Overview sheet: columns E and F (zh-cn / de-de status)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
